# Applies the "Added CDS All studies testcase" edit:
#  - Updates the SamplesTab query (cell B3) so it no longer selects the
#    Tumor / Analyte Type columns (sample_tumor_status, sample_type),
#    matching the trimmed-down "all studies" query from the diff.
#  - Moves the current selection / view down one row (was C2, now B3),
#    reflecting that the user was working on the SamplesTab query cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesTabQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND gi.instrument_model = 'Illumina HiSeq 2000'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newSamplesTabQuery

# Reflect the updated view/selection state captured in the workbook
$ws.Activate()
try {
    $excel.ActiveWindow.TopLeftCell = $ws.Range("A3")
} catch {
}
$ws.Range("B3").Select()
